$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 277.8902526399997
$ws.Range("E2").Value = 29092.72506141524
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 14865.25872276
$ws.Range("L2").Value = 50912.59821312752
$ws.Range("M2").Value = 11247.09127927
$ws.Range("N2").Value = 7270.39941619107
$ws.Range("O2").Value = 6890.515200515623

$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 5537.639684330663
$ws.Range("E2").Value = 56025.09626473462
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 37060.2472718608
$ws.Range("L2").Value = 73413.93799662068
$ws.Range("M2").Value = 21807.55922437125
$ws.Range("N2").Value = 10913.44862569963
$ws.Range("O2").Value = 9420.95505295181

$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 2266.487971660737
$ws.Range("B2").Value = 7112.344852924754
$ws.Range("E2").Value = 67064.77115738479
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 55915.31069905476
$ws.Range("L2").Value = 73413.93799662068
$ws.Range("M2").Value = 27623.86391009025
$ws.Range("N2").Value = 15939.31202052395
$ws.Range("O2").Value = 14199.58055466296

$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 2266.487971660737
$ws.Range("B2").Value = 7112.344852924754
$ws.Range("E2").Value = 67064.77115738479
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 55915.31069905476
$ws.Range("L2").Value = 73413.93799662068
$ws.Range("M2").Value = 27623.86391009025
$ws.Range("N2").Value = 15939.31202052395
$ws.Range("O2").Value = 16475.45838331085

$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 2266.487971660737
$ws.Range("B2").Value = 7112.344852924754
$ws.Range("E2").Value = 67064.77115738479
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 55915.31069905476
$ws.Range("L2").Value = 73413.93799662068
$ws.Range("M2").Value = 27623.86391009025
$ws.Range("N2").Value = 15939.31202052395
$ws.Range("O2").Value = 16475.45838331085

$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 2266.487971660737
$ws.Range("B2").Value = 7112.344852924754
$ws.Range("E2").Value = 67064.77115738479
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 55915.31069905476
$ws.Range("L2").Value = 73413.93799662068
$ws.Range("M2").Value = 27623.86391009025
$ws.Range("N2").Value = 15939.31202052395
$ws.Range("O2").Value = 16475.45838331085
